$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.194.69'
$ws.Range("E2").Value = '  -6.11%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.294.81'
$ws.Range("E3").Value = '  -5.13%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '558.89'
$ws.Range("E5").Value = '  -4.28%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.93'
$ws.Range("E6").Value = '  -2.40%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.296.11'
$ws.Range("E8").Value = '  -5.10%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.472'
$ws.Range("E9").Value = '  -1.96%  '

# Row 10
$ws.Range("E10").Value = '  -4.12%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.116'
$ws.Range("E11").Value = '  -5.14%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.370'
$ws.Range("E12").Value = '  -4.17%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.858.30'
$ws.Range("E13").Value = '  -5.02%  '

# Row 14
$ws.Range("E14").Value = '  -0.13%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.292.81'
$ws.Range("E15").Value = '  -5.14%  '

# Row 16
$ws.Range("E16").Value = '  -5.96%  '

# Row 17
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.377.91'
$ws.Range("E17").Value = '  -5.76%  '

# Row 18
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '24.10'
$ws.Range("E18").Value = '  -0.93%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.62'
$ws.Range("E19").Value = '  -0.88%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.30'
$ws.Range("E20").Value = '  -0.82%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.03'
$ws.Range("E21").Value = '  -9.37%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '350.79'
$ws.Range("E22").Value = '  -8.66%  '

# Row 23
$ws.Range("E23").Value = '  -2.80%  '

# Row 24
$ws.Range("E24").Value = '  -0.04%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.422.99'
$ws.Range("E25").Value = '  -5.24%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.22'
$ws.Range("E26").Value = '  -7.50%  '

# Row 27
$ws.Range("E27").Value = '  -3.29%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.09%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.17'
$ws.Range("E29").Value = '  +1.65%  '

# Row 30
$ws.Range("E30").Value = '  -0.43%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.78'
$ws.Range("E31").Value = '  -1.93%  '

# Row 32
$ws.Range("E32").Value = '  -6.05%  '

# Row 33
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.01%  '

# Row 34
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.149'
$ws.Range("E34").Value = '  -2.17%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.326.09'
$ws.Range("E35").Value = '  -4.97%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.66'
$ws.Range("E36").Value = '  -1.08%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.18'
$ws.Range("E37").Value = '  -0.21%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.72'
$ws.Range("E38").Value = '  -0.25%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.47'
$ws.Range("E39").Value = '  -1.64%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '158.43'
$ws.Range("E40").Value = '  -2.14%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0744'
$ws.Range("E41").Value = '  -4.31%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.10%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.01'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.34'
$ws.Range("E44").Value = '  +1.14%  '

# Row 46
$ws.Range("E46").Value = '  +0.92%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.84'
$ws.Range("E47").Value = '  -2.93%  '

# Row 48
$ws.Range("E48").Value = '  -4.79%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.64'
$ws.Range("E49").Value = '  -0.91%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.856'
$ws.Range("E50").Value = '  -5.00%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.29'
$ws.Range("E51").Value = '  +4.03%  '
